# Add a new automation-task row (Arti Gharjale) under the existing
# LoginData header/row on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Arti"
$ws.Range("B2").Value = "Gharjale"
# Leading apostrophe forces this numeric-looking value to be stored as
# text (quote-prefixed), matching the "431401" string in the sheet.
$ws.Range("C2").Value = "'431401"

# Leave the selection where the author left it.
$ws.Range("D2").Select() | Out-Null
